$wb = $excel.ActiveWorkbook

# -----------------------------------------------------------------------
# "Overview" sheet: row 3 corresponds to b.md.
#  E3 = zh-cn status, F3 = de-de status, G3 = Latest HO Xliff Generate Date
# -----------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E3").Value = "Ready for handoff"
$overview.Range("F3").Value = "Ready for handoff"
$overview.Range("G3").Value = "2016-08-19 06:36:57"

# -----------------------------------------------------------------------
# "zh-cn" sheet: row 3 (b.md) gets a fresh handoff.
# -----------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = "Ready for handoff"
# Leading apostrophe forces literal text so "False" isn't auto-promoted to a
# real Boolean; resetting the style afterwards keeps the cell's format the
# same as its neighbours (quotePrefix would otherwise stick around).
$zhcn.Range("F3").Value = "'False"
$zhcn.Range("F3").Style = "Normal"
$zhcn.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$zhcn.Range("H3").Value = "2016-08-19 06:36:52"
$zhcn.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c48e3eba6053ec07966460dcf971dea04e84d329/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f36674f61d39af4a2b26db496816facb1224eb8d/e2e/b.md."
$zhcn.Columns.Item(16).ColumnWidth = 39.166666666666664

# -----------------------------------------------------------------------
# "de-de" sheet: row 3 (b.md) gets a fresh handoff.
# -----------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = "Ready for handoff"
$dede.Range("F3").Value = "'False"
$dede.Range("F3").Style = "Normal"
$dede.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$dede.Range("H3").Value = "2016-08-19 06:36:57"
$dede.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c48e3eba6053ec07966460dcf971dea04e84d329/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f36674f61d39af4a2b26db496816facb1224eb8d/e2e/b.md."
$dede.Columns.Item(16).ColumnWidth = 39.166666666666664
